$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.364.11"
$ws.Range("E2").Value = "  -4.01%  "
$ws.Range("D3").Value = "2.360.03"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").Value = "  -0.05%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "511.03"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -4.30%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "128.36"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -4.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "2.378.23"
$ws.Range("E9").Value = "  -5.56%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.0955"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -3.99%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.152"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  -8.30%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.316"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -5.10%  "
$ws.Range("D14").Value = "2.782.92"
$ws.Range("E14").Value = "  -5.55%  "
$ws.Range("D15").Value = "56.256.74"
$ws.Range("E15").Value = "  -4.25%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "21.46"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -3.95%  "
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").Value = "2.359.91"
$ws.Range("E18").Value = "  -6.04%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "10.27"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -3.42%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "4.06"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -4.30%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "311.20"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "6.19"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "0.998"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "64.96"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -0.62%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.390"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -5.26%  "
$ws.Range("D27").Value = "2.465.71"
$ws.Range("E27").Value = "  -6.27%  "
$ws.Range("E28").Value = "  -4.87%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "7.22"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -3.73%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "174.53"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.68"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").Value = "0.0₃0714"
$ws.Range("E32").Value = "  -6.21%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "6.16"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("E34").Value = "  -7.52%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -2.77%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "1.20"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -4.03%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "3.70"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -6.39%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "35.64"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.790"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -6.35%  "
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "127.11"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "4.87"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -5.96%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "255.00"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -7.96%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.568"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -4.22%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.0901"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -3.52%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.0486"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("E50").Value = "  -6.10%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "16.57"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -6.53%  "
